# Applies the updated crypto price/volume snapshot to Sheet1.
# Several "Price" column values look numeric (e.g. "1.00", "0.0590") but must
# be preserved as literal text (matching the source inlineStr cells), so for
# those cells we temporarily force a text NumberFormat before assigning the
# value, then restore the default "Normal" style afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '63.367.48'
    'E2' = '  -3.38%  '
    'D3' = '2.603.49'
    'E3' = '  -1.67%  '
    'E4' = '  +0.17%  '
    'D5' = '571.90'
    'E5' = '  -3.98%  '
    'D6' = '154.11'
    'E6' = '  -1.16%  '
    'B7' = 'USDC'
    'C7' = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
    'D7' = '1.00'
    'E7' = '  +0.15%  '
    'B8' = 'XRP'
    'C8' = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
    'D8' = '0.630'
    'E8' = '  +0.27%  '
    'D9' = '0.120'
    'E9' = '  -5.53%  '
    'D10' = '5.74'
    'E10' = '  -1.28%  '
    'D11' = '0.383'
    'E11' = '  -3.28%  '
    'E12' = '  -0.53%  '
    'D13' = '28.09'
    'E13' = '  -2.08%  '
    'D14' = '3.083.68'
    'E14' = '  -1.33%  '
    'D15' = '0.0000182'
    'E15' = '  -7.28%  '
    'D16' = '63.337.54'
    'E16' = '  -3.18%  '
    'D17' = '2.610.74'
    'E17' = '  -1.33%  '
    'D18' = '12.03'
    'E18' = '  -4.17%  '
    'D19' = '4.59'
    'E19' = '  -2.76%  '
    'D20' = '7.43'
    'E20' = '  +0.06%  '
    'D21' = '340.88'
    'E21' = '  -2.12%  '
    'E22' = '  -0.14%  '
    'D23' = '66.91'
    'E23' = '  -2.97%  '
    'D24' = '1.73'
    'E24' = '  +2.44%  '
    'D25' = '0.0000107'
    'E25' = '  -4.24%  '
    'D26' = '9.20'
    'E26' = '  -4.48%  '
    'D27' = '565.16'
    'E27' = '  +7.00%  '
    'D28' = '1.55'
    'E28' = '  -2.16%  '
    'E29' = '  +0.20%  '
    'E30' = '  -2.69%  '
    'D31' = '7.78'
    'E31' = '  -1.14%  '
    'D32' = '2.05'
    'E32' = '  -3.17%  '
    'D33' = '1.69'
    'E33' = '  -3.40%  '
    'D34' = '6.40'
    'E34' = '  +0.18%  '
    'D35' = '5.26'
    'E35' = '  -2.77%  '
    'D36' = '0.406'
    'E36' = '  -3.38%  '
    'B37' = 'FirstDigitalUSD'
    'C37' = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    'D37' = '1.00'
    'E37' = '  +0.08%  '
    'B38' = 'EthereumClassic'
    'C38' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D38' = '19.73'
    'E38' = '  -3.08%  '
    'D39' = '151.56'
    'E39' = '  -3.45%  '
    'D40' = '1.85'
    'E40' = '  -3.65%  '
    'D41' = '1.00'
    'E41' = '  -0.01%  '
    'D42' = '41.64'
    'E42' = '  -1.76%  '
    'B43' = 'Aave'
    'C43' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D43' = '155.16'
    'E43' = '  -3.42%  '
    'B44' = 'dogwifhat'
    'C44' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'D44' = '2.34'
    'E44' = '  +3.52%  '
    'D45' = '3.92'
    'E45' = '  -3.43%  '
    'B46' = 'Hedera'
    'C46' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D46' = '0.0590'
    'E46' = '  -2.44%  '
    'B47' = 'InjectiveProtocol'
    'C47' = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    'D47' = '22.68'
    'E47' = '  +0.32%  '
    'D48' = '0.101'
    'E48' = '  +2.02%  '
    'D49' = '0.626'
    'E49' = '  -1.19%  '
    'D50' = '0.0248'
    'E50' = '  -2.07%  '
    'D51' = '18.84'
    'E51' = '  -4.52%  '
}

$forceTextCells = @(
    'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D13', 'D15', 'D18',
    'D19', 'D20', 'D21', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D31',
    'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41',
    'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51'
)

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    if ($forceTextCells -contains $addr) {
        $cell.NumberFormat = "@"
        $cell.Value = $updates[$addr]
        $cell.Style = "Normal"
    } else {
        $cell.Value = $updates[$addr]
    }
}
